# New "Pay Bill" test-data scenario: add a "Bills" worksheet (after "Login")
# with Company / Ref. No. / Bill ID headers and a sample SNGPL bill row, and
# refresh the stale "Cantt@mc.com" login fixture to "dha@mc.com".
#
# NOTE on write order: the shared-string table is built in first-write order,
# and the diff's target order (Abc@1234, Company, Ref. No., Bill ID, SNGPL,
# dha@mc.com, 37522010000) is only reproduced if we populate the new sheet's
# headers/first row before touching Login!A1 and before writing the bill id.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Login")

# Insert the new sheet right after "Login" and make it the active tab
# (mirrors Worksheets.Add placing + activating it, like a user adding a tab).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Bills"

# Header row
$ws2.Range("A1").Value = "Company"
$ws2.Range("B1").Value = "Ref. No."
$ws2.Range("C1").Value = "Bill ID"

# Sample data row
$ws2.Range("A2").Value = "SNGPL"

# Fix the stale login fixture on the "Login" sheet (A1: Cantt@mc.com -> dha@mc.com)
$ws1.Range("A1").Value = "dha@mc.com"

# Bill id is long enough that Excel would otherwise treat it as a number;
# the leading apostrophe forces text-with-quote-prefix storage, matching the
# new quotePrefix cell style in the target workbook.
$ws2.Range("B2").Value = "'37522010000"

# Column B needs to be wide enough to show the "Ref. No." header / bill id.
$ws2.Columns.Item(2).ColumnWidth = 18.5

# Leave the new sheet's selection on the bill-id cell.
$ws2.Range("B2").Select()
